# The three use-case sheets are renumbered down by one:
#   CU_10 -> CU_9
#   CU_11 -> CU_10
#   CU_12 -> CU_11
# Each sheet also carries its own case id as plain text in cell B2
# ("Caso de uso" row), so that label is updated to match the new
# sheet/tab name.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "CU_9"
$ws1.Range("B2").Value = "CU_9"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "CU_10"
$ws2.Range("B2").Value = "CU_10"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "CU_11"
$ws3.Range("B2").Value = "CU_11"
